{"js": "// Update the date line and the twenty-five \"NNN\u00d7N=\" multiplication\n// problems in the practice sheet to the new values from the target\n// revision. Every old string in the table is unique, so an exact\n// (case-sensitive, whole-match) search-and-replace on each pair safely\n// targets only the intended run.\nconst replacements = [\n  [\"2025-08-01 Friday\", \"2025-08-02 Saturday\"],\n  [\"227\u00d72=\", \"536\u00d76=\"],\n  [\"644\u00d75=\", \"406\u00d74=\"],\n  [\"371\u00d75=\", \"555\u00d76=\"],\n  [\"970\u00d79=\", \"440\u00d79=\"],\n  [\"475\u00d77=\", \"354\u00d77=\"],\n  [\"365\u00d74=\", \"472\u00d79=\"],\n  [\"740\u00d76=\", \"251\u00d73=\"],\n  [\"141\u00d77=\", \"317\u00d72=\"],\n  [\"870\u00d79=\", \"317\u00d75=\"],\n  [\"965\u00d79=\", \"254\u00d79=\"],\n  [\"882\u00d78=\", \"315\u00d75=\"],\n  [\"414\u00d72=\", \"402\u00d77=\"],\n  [\"785\u00d79=\", \"957\u00d74=\"],\n  [\"492\u00d79=\", \"907\u00d79=\"],\n  [\"383\u00d73=\", \"684\u00d76=\"],\n  [\"246\u00d74=\", \"436\u00d74=\"],\n  [\"186\u00d79=\", \"805\u00d79=\"],\n  [\"424\u00d79=\", \"458\u00d74=\"],\n  [\"405\u00d78=\", \"997\u00d73=\"],\n  [\"852\u00d75=\", \"388\u00d75=\"],\n  [\"229\u00d78=\", \"851\u00d73=\"],\n  [\"882\u00d76=\", \"639\u00d77=\"],\n  [\"809\u00d78=\", \"109\u00d75=\"],\n  [\"725\u00d76=\", \"528\u00d72=\"],\n  [\"356\u00d72=\", \"179\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# wdReplaceOne = 1 (replace only the current/first match found by this Execute call)\n$wdReplaceOne = 1\n\n$replacements = @(\n    @('2025-08-01 Friday', '2025-08-02 Saturday'),\n    @('227\u00d72=', '536\u00d76='),\n    @('644\u00d75=', '406\u00d74='),\n    @('371\u00d75=', '555\u00d76='),\n    @('970\u00d79=', '440\u00d79='),\n    @('475\u00d77=', '354\u00d77='),\n    @('365\u00d74=', '472\u00d79='),\n    @('740\u00d76=', '251\u00d73='),\n    @('141\u00d77=', '317\u00d72='),\n    @('870\u00d79=', '317\u00d75='),\n    @('965\u00d79=', '254\u00d79='),\n    @('882\u00d78=', '315\u00d75='),\n    @('414\u00d72=', '402\u00d77='),\n    @('785\u00d79=', '957\u00d74='),\n    @('492\u00d79=', '907\u00d79='),\n    @('383\u00d73=', '684\u00d76='),\n    @('246\u00d74=', '436\u00d74='),\n    @('186\u00d79=', '805\u00d79='),\n    @('424\u00d79=', '458\u00d74='),\n    @('405\u00d78=', '997\u00d73='),\n    @('852\u00d75=', '388\u00d75='),\n    @('229\u00d78=', '851\u00d73='),\n    @('882\u00d76=', '639\u00d77='),\n    @('809\u00d78=', '109\u00d75='),\n    @('725\u00d76=', '528\u00d72='),\n    @('356\u00d72=', '179\u00d73='),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceOne)\n}\n\nWrite-Output \"done\"\n"}
